$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8084924817085266
$ws.Range("B1").Value = 2.100004434585571
$ws.Range("D1").Value = 1.289513349533081
$ws.Range("E1").Value = 0.5090140104293823
